$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2311936337545805
$ws.Range("D2").Value = 0.819301027146599

$ws.Range("C3").Value = 0.7259912297960728
$ws.Range("D3").Value = 0.4754974572130677

$ws.Range("C4").Value = 1.845074836922812
$ws.Range("D4").Value = 0.07853162642023404
$ws.Range("G4").Value = "No"

$ws.Range("C5").Value = 2.94676437889888
$ws.Range("D5").Value = 0.007458250450645965

$ws.Range("C6").Value = 0.5429301262456993
$ws.Range("D6").Value = 0.5926365800552444

$ws.Range("C7").Value = 2.129855911196404
$ws.Range("D7").Value = 0.04461406075880214

$ws.Range("C8").Value = 3.173636740993039
$ws.Range("D8").Value = 0.004396590026047198

$ws.Range("C9").Value = 1.228911409267549
$ws.Range("D9").Value = 0.2320884271663823

$ws.Range("C10").Value = 2.947826918743542
$ws.Range("D10").Value = 0.007439998481404198

$ws.Range("C11").Value = 0.5955171939783672
$ws.Range("D11").Value = 0.557574967669715
